# Gray out the "File I/O" section's sub-heading and its two bullet points,
# matching the existing "dimmed" style already used elsewhere in the doc
# (color A6A6A6 / theme color background1 / theme shade A6), applied to
# both the paragraph mark (pPr/rPr) and the run (r/rPr).

$d = $word.ActiveDocument

$colorXml = '<w:rPr><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/></w:rPr>'

function Apply-DimColor([string]$searchText) {
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "Not found: $searchText"
        return
    }

    # Work with the whole paragraph range (includes the end-of-paragraph
    # mark) so we can stamp both the paragraph-mark run properties and the
    # text run's properties in one shot.
    $para = $rng.Paragraphs(1)
    $prng = $para.Range
    $openXml = $prng.WordOpenXML

    if (-not ($openXml -match '(?s)(<w:p[ >].*?</w:p>)')) {
        Write-Host "Could not isolate paragraph XML for: $searchText"
        return
    }
    $pxml = $matches[1]

    # Add/extend <w:pPr> with the color rPr (paragraph mark formatting).
    if ($pxml -match '(?s)<w:pPr\s*/>') {
        $pxml = $pxml -replace '(?s)<w:pPr\s*/>', ('<w:pPr>' + $colorXml + '</w:pPr>')
    } elseif ($pxml -match '(?s)<w:pPr>.*?</w:pPr>') {
        $pxml = $pxml -replace '(?s)(<w:pPr>.*?)</w:pPr>', ('$1' + $colorXml + '</w:pPr>')
    } else {
        $pxml = $pxml -replace '(?s)^(<w:p(?:\s[^>]*)?>)', ('$1<w:pPr>' + $colorXml + '</w:pPr>')
    }

    # Add the same color rPr as the first child of every run in the
    # paragraph that doesn't already carry its own rPr (covers both bare
    # <w:r> and <w:r someAttr="..."> forms).
    $pxml = $pxml -replace '(?s)(<w:r(?:\s[^>]*)?>)(?!\s*<w:rPr>)', ('$1' + $colorXml)

    $wrapper = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData><w:document ' +
               'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
               'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
               'xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml">' +
               '<w:body>' + $pxml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $prng.Text = ""
    [void]$prng.InsertXML($wrapper)
}

Apply-DimColor("Lưu dữ liệu vào File:")
Apply-DimColor("Số tiền của Nhân vật")
Apply-DimColor("Những hạt giống mà nhân vật có")
